# Apply updated cryptocurrency price/volume data to sheet1
# (diff: Updated cryptos list with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.609.07"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.960.56"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.77"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.93"
$ws.Range("E6").Value = "  -0.51%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.961.45"
$ws.Range("E8").Value = "  +1.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  +4.99%  "

# Row 11
$ws.Range("E11").Value = "  -0.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  +1.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  +3.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.43"
$ws.Range("E14").Value = "  -1.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.126"
$ws.Range("E15").Value = "  -0.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.453.00"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.506.42"
$ws.Range("E17").Value = "  +2.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.72"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.955.55"
$ws.Range("E19").Value = "  +0.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "440.98"
$ws.Range("E20").Value = "  +0.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.47"
$ws.Range("E21").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.674"
$ws.Range("E22").Value = "  -0.85%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.11"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.79"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.09"
$ws.Range("E25").Value = "  +0.78%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.93"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("E27").Value = "  -3.36%  "

# Row 28
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("E29").Value = "  +0.13%  "

# Row 30
$ws.Range("E30").Value = "  +0.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.13"
$ws.Range("E31").Value = "  -6.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.42"
$ws.Range("E32").Value = "  -0.93%  "

# Row 33
$ws.Range("E33").Value = "  -1.86%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0876"
$ws.Range("E35").Value = "  +0.55%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -2.45%  "

# Row 37
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.04"
$ws.Range("E38").Value = "  +2.30%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.65"
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
$ws.Range("E40").Value = "  -2.95%  "

# Row 41
$ws.Range("E41").Value = "  -0.18%  "

# Row 42
$ws.Range("E42").Value = "  -5.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.283"
$ws.Range("E43").Value = "  -1.96%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.42"
$ws.Range("E44").Value = "  -6.62%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.717.35"
$ws.Range("E45").Value = "  +0.78%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.07"
$ws.Range("E46").Value = "  +1.50%  "

# Row 47
$ws.Range("E47").Value = "  -1.64%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "366.43"
$ws.Range("E48").Value = "  -2.63%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.12"
$ws.Range("E50").Value = "  -3.76%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  -0.47%  "
